# Apply building.xlsx update:
#  - params sheet: effective_heat_capacity 200 -> 60, differential_cost 0 -> 75,
#    add new row 7 "WRG" = 75 %
#  - thermal_hull sheet: update U-Wert for Wand, restructure the Boden/Dach/Fenster rows
#    (insert a new row so the table grows from 5 to 6 data rows) and add new values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: params
# ---------------------------------------------------------------------------
$params = $wb.Worksheets.Item("params")

$params.Range("B4").Value = 60
$params.Range("B6").Value = 75

$params.Range("A7").Value = "WRG"
$params.Range("B7").Value = 75
$params.Range("C7").Value = "%"

$params.Activate()
$params.Range("A8").Select()

# ---------------------------------------------------------------------------
# Sheet 2: thermal_hull
# ---------------------------------------------------------------------------
$hull = $wb.Worksheets.Item("thermal_hull")

# U-Wert of Wand changes from 0.3 to 0.6
$hull.Range("C2").Value = 0.6

# Insert a new row above row 3 so Boden/Dach/Fenster all move down one row,
# then rebuild the Boden/Dach/Fenster rows with their new layout/values.
$hull.Rows("3").Insert()

$hull.Range("A3").Value = 1
$hull.Range("B3").Value = "Boden"
$hull.Range("C3").Value = 800
$hull.Range("D3").Value = 1

$hull.Range("A4").Value = 0.3
$hull.Range("B4").Value = 1
$hull.Range("C4").Value = "Dach"
$hull.Range("D4").Value = 1

$hull.Range("A5").Value = 800
$hull.Range("B5").Value = 0.6
$hull.Range("C5").Value = 1
$hull.Range("D5").Value = 1

$hull.Range("A6").Value = "Fenster"
$hull.Range("B6").Value = 0
$hull.Range("C6").Value = 0
$hull.Range("D6").ClearContents()

$hull.Activate()
$hull.Range("A11").Select()
